# Auto-upload VRF Excel file: add a new "test" sheet with sample VRF model data.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "test"

# Header row - same labels used by every other sheet in this workbook.
$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"

# Bold, bordered, centered header formatting (matches the other sheets' header style).
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Sample outdoor-unit row.
$ws.Range("A2").Value = "asd"
$ws.Range("B2").Value = 1
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "11"

# Sample indoor-unit row.
$ws.Range("D3").Value = "asdf"
$ws.Range("E3").Value = 1
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "12"

$ws.Range("A1").Select()
